# Apply the benchmark-table corrections described by the commit.
# The document contains a single table, one column, one value per row.
# Word COM collections are 1-indexed.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Simple single-value row updates (row -> new text)
$updates = @{
    1  = "0M"
    2  = "0M"
    3  = "0M"
    4  = "1922"
    5  = "0.00001"
    6  = "0.00121"
    7  = "0.00016"
    8  = "0.00005"
    9  = "0.00027"
    10 = "0.00032"
    11 = "0.00041"
    12 = "0.36399"
}

foreach ($rowIndex in $updates.Keys) {
    $cell = $t.Cell($rowIndex, 1)
    $cell.Range.Text = $updates[$rowIndex]
}

# Rows 44-46 previously held tab-separated multi-value summaries; they are
# collapsed down to a single value each (matching rows 1-3's originals).
$t.Cell(44, 1).Range.Text = "99.74"
$t.Cell(45, 1).Range.Text = "0.36"
$t.Cell(46, 1).Range.Text = "142"
